$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = " Здравствуйте. Пожалуйста, свяжитесь с нами в чате поддержки по WhatsApp 7 (962) 559-29-48. Проверим товар, если брак подтвердиться мы сделаем возврат"
$ws.Range("A2").Value = "Нет заказа "
$ws.Range("B2").Value = "Здравствуйте. мы отгружаем товар со своего склада на Вб и их логистика делает доставку. напиши в тех.поддержку самого Вб, мы не можем повлиять на скорость доставки, ни отменить ее не можем. отслеживайте его доставку в своем личном кабинете"

$ws.Range("A1:B2").WrapText = $true
$ws.Rows.Item(1).RowHeight = 57.6
$ws.Rows.Item(2).RowHeight = 86.4

$ws.Range("A2:B2").Select()
